{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block,\n// along with the blank paragraph that separates it from the bibliography\n// entry above it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the SERAFINI bibliography entry (last bibliography line) so we can\n// anchor the deletion to the blank paragraph that immediately follows it.\nlet serafiniIndex = -1;\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.indexOf(\"SERAFINI, Maria Jos\u00e9\") !== -1) {\n    serafiniIndex = i;\n  }\n  if (t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n  }\n  if (t.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (jupiterIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\"Could not locate footer paragraphs to remove.\");\n}\n\n// The blank paragraph right before \"Ver no Jupiter ...\" (and right after\n// the SERAFINI entry) is removed too.\nlet blankIndex = jupiterIndex - 1;\nif (serafiniIndex !== -1 && blankIndex === serafiniIndex) {\n  // no blank paragraph found where expected; nothing to delete for it\n  blankIndex = -1;\n} else if (blankIndex >= 0 && items[blankIndex].text !== \"\") {\n  blankIndex = -1;\n}\n\nconst toDelete = [];\nif (blankIndex !== -1) toDelete.push(blankIndex);\ntoDelete.push(jupiterIndex);\ntoDelete.push(copyrightIndex);\n\n// Delete from the highest index to the lowest so earlier indices remain\n// valid while we work.\ntoDelete.sort((a, b) => b - a);\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"(c) 2020 ...\" footer block,\n# along with the blank paragraph that separates it from the bibliography\n# entry (\"SERAFINI, ...\") above it.\n\n$d = $word.ActiveDocument\n\n$serafiniIndex = -1\n$jupiterIndex = -1\n$copyrightIndex = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith(\"SERAFINI, Maria\")) {\n        $serafiniIndex = $i\n    }\n    elseif ($t.StartsWith(\"Ver no Jupiter\")) {\n        $jupiterIndex = $i\n    }\n    elseif ($t.Contains(\"Contact: luizeleno@usp.br\")) {\n        $copyrightIndex = $i\n    }\n}\n\nif ($jupiterIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"Could not locate footer paragraphs to remove.\"\n}\n\n$blankIndex = $jupiterIndex - 1\nif ($serafiniIndex -ne -1 -and $blankIndex -eq $serafiniIndex) {\n    $blankIndex = -1\n}\nelseif ($blankIndex -ge 1) {\n    $blankText = $d.Paragraphs.Item($blankIndex).Range.Text\n    if ($blankText -ne [string][char]13) {\n        $blankIndex = -1\n    }\n}\n\n$toDelete = @($copyrightIndex, $jupiterIndex)\nif ($blankIndex -ne -1) {\n    $toDelete += $blankIndex\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$toDelete = $toDelete | Sort-Object -Descending\nforeach ($idx in $toDelete) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
